$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 176, shifting existing rows 176-279 down to 179-282
$ws.Rows("176:178").Insert()

# Row 176 (new record)
$ws.Range("A176").Value = 5
$ws.Range("B176").Value = "Macroferia Regional de Talca"
$ws.Range("C176").Value = "Maule"
$ws.Range("D176").Value = 44455
$ws.Range("E176").Value = 7
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100101
$ws.Range("H176").Value = "Berries"
$ws.Range("I176").Value = 100112025
$ws.Range("J176").Value = "Frutilla"
$ws.Range("K176").Value = "Sin especificar"
$ws.Range("L176").Value = "Especial"
$ws.Range("M176").Value = 160
$ws.Range("N176").Value = 23000
$ws.Range("O176").Value = 24000
$ws.Range("P176").Value = 23625
$ws.Range("Q176").Value = "`$/bandeja 7 kilos"
$ws.Range("R176").Value = "Provincia de Melipilla"
$ws.Range("S176").Value = 3375
$ws.Range("T176").Value = 7

# Row 177 (new record)
$ws.Range("A177").Value = 5
$ws.Range("B177").Value = "Macroferia Regional de Talca"
$ws.Range("C177").Value = "Maule"
$ws.Range("D177").Value = 44455
$ws.Range("E177").Value = 7
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100101
$ws.Range("H177").Value = "Berries"
$ws.Range("I177").Value = 100112025
$ws.Range("J177").Value = "Frutilla"
$ws.Range("K177").Value = "Sin especificar"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 100
$ws.Range("N177").Value = 20000
$ws.Range("O177").Value = 20000
$ws.Range("P177").Value = 20000
$ws.Range("Q177").Value = "`$/bandeja 7 kilos"
$ws.Range("R177").Value = "Provincia de Melipilla"
$ws.Range("S177").Value = 2857
$ws.Range("T177").Value = 7

# Row 178 (new record)
$ws.Range("A178").Value = 5
$ws.Range("B178").Value = "Macroferia Regional de Talca"
$ws.Range("C178").Value = "Maule"
$ws.Range("D178").Value = 44455
$ws.Range("E178").Value = 7
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100101
$ws.Range("H178").Value = "Berries"
$ws.Range("I178").Value = 100112025
$ws.Range("J178").Value = "Frutilla"
$ws.Range("K178").Value = "Sin especificar"
$ws.Range("L178").Value = "Segunda"
$ws.Range("M178").Value = 60
$ws.Range("N178").Value = 12000
$ws.Range("O178").Value = 12000
$ws.Range("P178").Value = 12000
$ws.Range("Q178").Value = "`$/bandeja 7 kilos"
$ws.Range("R178").Value = "Provincia de Melipilla"
$ws.Range("S178").Value = 1714
$ws.Range("T178").Value = 7
